# Applies the updated Transavia daily production dataset (UTC -> EET fix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    45506, 45506, 45506, 45506, 45506, 45506, 45506, 45506, 45506, 45506,
    45506, 45506, 45506, 45506, 45507, 45507, 45507, 45507, 45507, 45507,
    45507, 45507, 45507, 45507, 45507, 45507, 45507, 45507, 45507, 45507,
    45507, 45507, 45507, 45507, 45507, 45507, 45507, 45507, 45508, 45508,
    45508, 45508, 45508, 45508, 45508, 45508, 45508, 45508, 45508, 45508,
    45508, 45508, 45508, 45508, 45508, 45508, 45508, 45508, 45508, 45508,
    45508, 45508, 45509, 45509, 45509, 45509, 45509, 45509, 45509, 45509,
    45509, 45509, 45509, 45509, 45509, 45509, 45509, 45509, 45509, 45509,
    45509, 45509, 45509, 45509, 45509, 45509, 45510, 45510, 45510, 45510,
    45510, 45510, 45510, 45510, 45510, 45510, 45510, 45510, 45510, 45510,
    45510, 45510, 45510, 45510, 45510, 45510, 45510, 45510, 45510, 45510,
    45511, 45511, 45511, 45511, 45511, 45511, 45511, 45511, 45511, 45511,
    45511, 45511, 45511, 45511, 45511, 45511, 45511, 45511, 45511, 45511,
    45511, 45511, 45511, 45511, 45512, 45512, 45512, 45512, 45512, 45512,
    45512, 45512, 45512, 45512, 45512, 45512, 45512, 45512, 45512, 45512,
    45512, 45512, 45512, 45512, 45512, 45512, 45512, 45512, 45513, 45513,
    45513, 45513, 45513, 45513, 45513, 45513, 45513, 45513, 45513
)

$hours = @(
    10, 11, 12, 13, 14, 15, 16, 17, 18, 19,
    20, 21, 22, 23, 0, 1, 2, 3, 4, 5,
    6, 7, 8, 9, 10, 11, 12, 13, 14, 15,
    16, 17, 18, 19, 20, 21, 22, 23, 0, 1,
    2, 3, 4, 5, 6, 7, 8, 9, 10, 11,
    12, 13, 14, 15, 16, 17, 18, 19, 20, 21,
    22, 23, 0, 1, 2, 3, 4, 5, 6, 7,
    8, 9, 10, 11, 12, 13, 14, 15, 16, 17,
    18, 19, 20, 21, 22, 23, 0, 1, 2, 3,
    4, 5, 6, 7, 8, 9, 10, 11, 12, 13,
    14, 15, 16, 17, 18, 19, 20, 21, 22, 23,
    0, 1, 2, 3, 4, 5, 6, 7, 8, 9,
    10, 11, 12, 13, 14, 15, 16, 17, 18, 19,
    20, 21, 22, 23, 0, 1, 2, 3, 4, 5,
    6, 7, 8, 9, 10, 11, 12, 13, 14, 15,
    16, 17, 18, 19, 20, 21, 22, 23, 0, 1,
    2, 3, 4, 5, 6, 7, 8, 9, 10
)

$vals = @(
    0.01420320570468903, 0.0231179092079401, 0.7500006556510925, 0.7905964851379395, 0.7702441811561584, 0.6375821828842163,
    0.4212532341480255, 0.4932942688465118, 0.2481384724378586, 0.1314605325460434, 0.001458370708860457, -0.001172069576568902,
    -0.001172069576568902, -0.00104915932752192, -0.0001568608568049967, -0.0001568608568049967, -0.0001568608568049967, -0.0001568608568049967,
    -0.0001568608568049967, -0.0000711647080606781, 0.06904695928096771, 0.2204293459653854, 0.4437849521636963, 0.6803387999534607,
    0.8479066491127014, 0.8736859560012817, 0.8839192390441895, 0.8703525066375732, 0.7820696830749512, 0.760837197303772,
    0.5008438229560852, 0.4155340194702148, 0.3269176185131073, 0.1648519039154053, 0.0360427163541317, 0.0007388877565972507,
    -0.00104915932752192, -0.001010326785035431, -0.0001568608568049967, -0.0000711647080606781, 0.0002246848307549953, 0.0004331798991188407,
    0.0004688547342084348, 0.0004688547342084348, 0.1022904962301254, 0.3387464284896851, 0.3976575434207916, 0.7315188050270081,
    0.9679725766181946, 1.06139874458313, 0.9730522036552429, 0.7821691036224365, 0.6720830202102661, 0.3899971842765808,
    0.4613310992717743, 0.5470319390296936, 0.3222713470458984, 0.1415599137544632, 0.03071128018200397, 0.0007388877565972507,
    -0.001172069576568902, -0.00104915932752192, -0.0001568608568049967, -0.0001568608568049967, -0.0001568608568049967, -0.0001568608568049967,
    -0.0001568608568049967, -0.0001568608568049967, 0.09593164175748825, 0.2769536972045898, 0.4733706414699554, 0.7531857490539551,
    0.7379154562950134, 0.8077488541603088, 1.010447263717651, 1.01949155330658, 0.9143903851509094, 0.626031756401062,
    0.4507476687431335, 0.3424841463565826, 0.1986044347286224, 0.1432812511920929, 0.01688057743012905, -0.001172069576568902,
    -0.001172069576568902, -0.001010326785035431, -0.0001568608568049967, -0.0000711647080606781, -0.0000711647080606781, -0.0000711647080606781,
    -0.0000711647080606781, 0.0002246848307549953, 0.09664736688137054, 0.3484517931938171, 0.3808132410049438, 0.7315188050270081,
    0.9811878204345703, 1.06139874458313, 0.9730522036552429, 0.9585938453674316, 0.7264062166213989, 0.3772031962871552,
    0.3171989619731903, 0.4267469644546509, 0.2833408415317535, 0.1549337357282639, 0.02294223569333553, -0.001172069576568902,
    -0.001010326785035431, -0.001010326785035431, -0.0001568608568049967, -0.0001568608568049967, -0.0000711647080606781, -0.0000711647080606781,
    0.0002246848307549953, 0.0002246848307549953, 0.09614302217960358, 0.3305690586566925, 0.485541433095932, 0.7473618388175964,
    0.9136647582054138, 0.8964380621910095, 0.9864011406898499, 0.9816029071807861, 0.6205585598945618, 0.5767087936401367,
    0.3827455937862396, 0.4271636307239532, 0.3884811103343964, 0.1505521088838577, 0.03530823439359665, 0.01442866679280996,
    -0.0001870296255219728, -0.002097988268360496, -0.001460882951505482, -0.001460882951505482, -0.001337968744337559, -0.001337968744337559,
    -0.001337968744337559, -0.001299136434681714, 0.06688468158245087, 0.2543617486953735, 0.4838624894618988, 0.7057944536209106,
    0.847597062587738, 0.7063229084014893, 0.5931406617164612, 0.7142804265022278, 0.3891435265541077, 0.4497986733913422,
    0.3740514814853668, 0.4271636307239532, 0.2451140731573105, 0.2427791953086853, 0.02982072904706001, 0.0142905842512846,
    0.01535458210855722, 0.01535458210855722, 0.01700847409665585, 0.01700847409665585, 0.001592351705767214, 0.0004500749055296183,
    -0.001460882951505482, -0.001460882951505482, 0.08554156124591827, 0.4034946858882904, 0.5001636147499084, 0.8065139055252075,
    0.8376703262329102
)

$startRow = 2
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $dates[$i]
    $ws.Cells.Item($r, 2).Value2 = $hours[$i]
    $ws.Cells.Item($r, 3).Value2 = $vals[$i]
}

Write-Host "Updated $($dates.Length) rows."
